{"js": "// Update the East-Asian / Complex-Script font fallbacks recorded in the\n// document's paragraph styles (vignettes/docx/text-font-size.docx):\n//   - \"Normal\" and \"Heading\": eastAsia font DejaVu Sans -> Tahoma\n//   - \"List\", \"Caption\", \"Index\": pick up an explicit complex-script (cs)\n//     font of \"DejaVu Sans\" (previously falling back to the doc default)\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal\");\nawait context.sync();\n\nfunction getStyle(name) {\n  const match = styles.items.find((s) => s.nameLocal === name);\n  if (!match) {\n    throw new Error(\"style not found: \" + name);\n  }\n  return match;\n}\n\nconst normalStyle = getStyle(\"Normal\");\nconst headingStyle = getStyle(\"Heading\");\nconst listStyle = getStyle(\"List\");\nconst captionStyle = getStyle(\"Caption\");\nconst indexStyle = getStyle(\"Index\");\n\nnormalStyle.font.nameFarEast = \"Tahoma\";\nheadingStyle.font.nameFarEast = \"Tahoma\";\n\nlistStyle.font.nameBidirectional = \"DejaVu Sans\";\ncaptionStyle.font.nameBidirectional = \"DejaVu Sans\";\nindexStyle.font.nameBidirectional = \"DejaVu Sans\";\n\nawait context.sync();\n", "ps1": "# Update the East-Asian / Complex-Script font fallbacks recorded in the\n# document's paragraph styles (vignettes/docx/text-font-size.docx):\n#   - \"Normal\" and \"Heading\": eastAsia font DejaVu Sans -> Tahoma\n#   - \"List\", \"Caption\", \"Index\": pick up an explicit complex-script (cs)\n#     font of \"DejaVu Sans\" (previously falling back to the doc default)\n$d = $word.ActiveDocument\n\n$normal = $d.Styles.Item(\"Normal\")\n$normal.Font.NameFarEast = \"Tahoma\"\n\n$heading = $d.Styles.Item(\"Heading\")\n$heading.Font.NameFarEast = \"Tahoma\"\n\n$list = $d.Styles.Item(\"List\")\n$list.Font.NameBi = \"DejaVu Sans\"\n\n$caption = $d.Styles.Item(\"Caption\")\n$caption.Font.NameBi = \"DejaVu Sans\"\n\n$index = $d.Styles.Item(\"Index\")\n$index.Font.NameBi = \"DejaVu Sans\"\n"}
